# Regenerate penyata (5LUHUR-2023) to follow new data and format.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Text relabeling (shared-string content updates)
# ---------------------------------------------------------------------------
# "Kali Pertama/Kedua/Ketiga/Keempat" -> "Semakan Kali ..." (used under three
# sections: Penandaan Fail, Laporan Atas Talian, JPPM / JDM / JDRM)
foreach ($r in 16, 22, 28) { $ws.Range("C$r").Value = "Semakan Kali Pertama" }
foreach ($r in 17, 23, 29) { $ws.Range("C$r").Value = "Semakan Kali Kedua" }
foreach ($r in 18, 24, 30) { $ws.Range("C$r").Value = "Semakan Kali Ketiga" }
foreach ($r in 19, 25, 31) { $ws.Range("C$r").Value = "Semakan Kali Keempat" }

# Competition names switched from all-caps to title case.
$ws.Range("C34").Value = "Bouquet Kreatif"
$ws.Range("C35").Value = "Tik Tok Raya"
$ws.Range("C36").Value = "Riang Ria Kuih Raya"
$ws.Range("C37").Value = "Creative Collage"

# ---------------------------------------------------------------------------
# 2. Updated transaction figures
# ---------------------------------------------------------------------------
$ws.Range("D18").Value = 7075
$ws.Range("E18").Value = 10705
$ws.Range("D23").Value = 400
$ws.Range("E29").Value = 750

# ---------------------------------------------------------------------------
# 3. Header block: title moves from E4 to D4 and now spans D4:G4
# ---------------------------------------------------------------------------
$title = $ws.Range("E4").Value2
$ws.Range("E4").ClearContents()
$ws.Range("D4").Value = $title

# ---------------------------------------------------------------------------
# 4. Merged-range layout changes
# ---------------------------------------------------------------------------
$ws.Range("D4:G4").Merge()
$ws.Range("B5:C5").Merge()
$ws.Range("B12:F12").Merge()

$ws.Range("B15:C15").UnMerge()

$ws.Range("B21:C21").UnMerge()
$ws.Range("B21:E21").Merge()
$ws.Range("B27:C27").UnMerge()
$ws.Range("B27:E27").Merge()
$ws.Range("B33:C33").UnMerge()
$ws.Range("B33:E33").Merge()

$ws.Range("B43:E43").Merge()

# ---------------------------------------------------------------------------
# 5. Page / print setup
# ---------------------------------------------------------------------------
$ws.PageSetup.FitToPagesTall = 1
$ws.PageSetup.FitToPagesWide = 1
$ws.PageSetup.CenterHorizontally = $true
$ws.PageSetup.HeaderMargin = 0
$ws.PageSetup.FooterMargin = 0

# ---------------------------------------------------------------------------
# 6. Picture (logo) repositioned / resized
# ---------------------------------------------------------------------------
$shp = $ws.Shapes.Item(1)
$shp.Left = 41.2125
$shp.Top = 14.25
$shp.Width = 46.5
$shp.Height = 47.25

# ---------------------------------------------------------------------------
# 7. Trim trailing blank formatted row
# ---------------------------------------------------------------------------
$ws.Rows.Item(1001).Delete()
